# Auto-generated PowerShell COM-interop script
$wb = $excel.ActiveWorkbook
$wsMeasures = $wb.Worksheets.Item("researchMeasures")
$wsNutrition = $wb.Worksheets.Item("NutritionalData")

# --- NutritionalData: new recipe rows 148-152 ---------------------------------
$wsNutrition.Range("A148").Value = 'spaghetti 3-11-2021 recipe 2 pkg Barilla Gluten free spaghetti serves 6, 1 jar Prego 3 cheese serves 5, 1 pkg beyond meat 4.5 servings, 2 bell peppers 1 yellow other green, 2 cups mushrooms, 2tbs olive oil, this makes about 4-5 bowls. I will say it makes 5 bowls, so that each bowl is 1/5th this pot of spaghetti.'
$wsNutrition.Range("B148").Formula = "=SUM(B49*2,B39*2,B36*4.5,B35*6,B114*5,B149*15)"
$wsNutrition.Range("C148").Formula = "=SUM(C49*2,C39*2,C36*4.5,C35*6,C114*5,C149*15)"
$wsNutrition.Range("D148").Formula = "=SUM(D49*2,D39*2,D36*4.5,D35*6,D114*5,D149*15)"
$wsNutrition.Range("E148").Formula = "=SUM(E49*2,E39*2,E36*4.5,E35*6,E114*5,E149*15)"
$wsNutrition.Range("F148").Formula = "=SUM(F49*2,F39*2,F36*4.5,F35*6,F114*5,F149*15)"
$wsNutrition.Range("G148").Formula = "=SUM(G49*2,G39*2,G36*4.5,G35*6,G114*5,G149*15)"
$wsNutrition.Range("H148").Formula = "=SUM(H49*2,H39*2,H36*4.5,H35*6,H114*5,H149*15)"

$wsNutrition.Range("A149").Value = 'white sliced mushrooms, calorieking.com, for 1 mushroom'
$wsNutrition.Range("B149").Value = 2
$wsNutrition.Range("C149").Value = 0.1
$wsNutrition.Range("D149").Value = 0.1
$wsNutrition.Range("E149").Value = 0.3
$wsNutrition.Range("F149").Value = 0.2
$wsNutrition.Range("G149").Value = 0.1
$wsNutrition.Range("H149").Value = 1

$wsNutrition.Range("A150").Value = 'bowl of spaghetti 3-11-2021 recipe divided into 5 bowls'
$wsNutrition.Range("B150").Formula = "=B148/5"
$wsNutrition.Range("C150").Formula = "=C148/5"
$wsNutrition.Range("D150").Formula = "=D148/5"
$wsNutrition.Range("E150").Formula = "=E148/5"
$wsNutrition.Range("F150").Formula = "=F148/5"
$wsNutrition.Range("G150").Formula = "=G148/5"
$wsNutrition.Range("H150").Formula = "=H148/5"

$wsNutrition.Range("A151").Value = 'mozzarella parmesan shredded cheese Good & Gather Target brand, serving sz 1/4 cup:'
$wsNutrition.Range("B151").Value = 100
$wsNutrition.Range("C151").Value = 6
$wsNutrition.Range("D151").Value = 4
$wsNutrition.Range("E151").Value = 8
$wsNutrition.Range("F151").Value = 2
$wsNutrition.Range("G151").Value = 0
$wsNutrition.Range("H151").Value = 280

$wsNutrition.Range("B152").Formula = "=B144*8/3"
$wsNutrition.Range("C152").Formula = "=C144*8/3"
$wsNutrition.Range("D152").Formula = "=D144*8/3"
$wsNutrition.Range("E152").Formula = "=E144*8/3"
$wsNutrition.Range("F152").Formula = "=F144*8/3"
$wsNutrition.Range("G152").Formula = "=G144*8/3"
$wsNutrition.Range("H152").Formula = "=H144*8/3"
$wsNutrition.Range("B152:H152").NumberFormat = "0.00"

# column A of these new rows mirrors the left/top aligned label style used below it
$wsNutrition.Range("A148:A151").HorizontalAlignment = -4131
$wsNutrition.Range("A148:A151").VerticalAlignment = -4160

# --- researchMeasures: update row 57 diary text + recomputed totals ----------
$wsMeasures.Range("AA57").Value = '3 eggs
(210	15	4.5	18	0	0	210)
2 tbs sourcream
(60	5	3.5	1	2	0	15)
2 tbs olive oil
(240	28	4	0	0	0	0)
1/3 green bell pepper
(13.3	0	0	0.3	3.3	1	0)
1/2 cup mozzarella cheese
(160	10	7	12	2	0	380)
8 mini corn tortillas
(266.67	 2.67	0.00	5.33	53.33	5.33	53.33)
2 yellow cake cupcakes
(270	12	2.6	3	36	0	310)
1 serving Tostitos hint of lime tortilla chips
(150	7	1	2	18	1	130)
2 tbs sourcream
(60	5	3.5	1	2	0	15)
1/4 cup mozzarella cheese
(80	5	3.5	6	1	0	190)
1 serving pea protein
(130	2	0	18	9	2	320)
1 banana
(105	0	0	1	27	3	1)
2 tbs cocoa
(20	1	0	2	6	2	0)
1 tbs honey
(60	0	0	0	17	0	0)
1 serving walnuts
(200	20	2	5	4	2	0)
1/4 cup Silk Amond milk unsweet
(7.5	0.625	0	0.25	0.25	0	28.75)
8 corn tortillas
(266.67	2.67	0.00	5.33	53.33	5.33	53.33)
1/3 cup mozzarella cheese
(80	5	3.5	6	1	0	190)
bowl of spaghetti made 3-11-2021
(614	24.8	6.1	26.1	72.9	5.5	798)
1/4 cup mozzarella Target brand
(100	6	4	8	2	0	280)
2 servings tostitos chips
(300	14	2	4	36	2	260)
5 servings sourcream about 2 tbs
(30	25	17.5	5	10	0	75)
=210+60+240+13.3+160+267+270+150+60+80+130+105+20+60+200+7.5+266.67+80+614+100+300+30
=15+5+28+0+10+2.67+12+7+5+5+2+0+1+0+20+0.625+2.67+5+24.8+6+14+25
=4.5+3.5+4+0+7+0+2.6+1+3.5+3.5+0+0+0+0+2+0+0+3.5+6.1+4+2+17.5
=18+1+0+0.3+12+5.33+3+2+1+6+18+1+2+0+5+0.25+5.33+6+26.1+8+4+5
=0+2+0+3.3+2+53.33+36+18+2+1+9+27+6+17+4+0.25+53.33+1+72.9+2+36+10
=0+0+0+1+0+5.33+0+1+0+0+2+3+2+0+2+0+5.33+0+5.5+0+2+0
=210+15+0+0+380+53.33+310+130+15+190+320+1+0+0+0+28.75+53.33+190+798+280+260+75
'
$wsMeasures.Range("Z57").Value = 'Woke up at 5:20 am, got out of bed at 545 am, made coffee and fed the babies after cleaning their messes first thing. Started a few of the remaining chemistry problems on moles produced from combustion and balancing equations, made another cup of coffee for the 2nd cup, and roommate arrived a little before 7 am. Finished my coffee and had a BM reg sz by 720 am, took a break from the last equation that was a long one, and missed 1st of 3 attempts on it. Had to balance an equation of propane given the density of propane in g/mL and find the mass of CO2 in kg as product in combustion after solving mass of propane and converting mole ratios of equation not given but assumed to be C3H3+O2-->CO2+H2O in some balanced coefficient assortment that somewhere and something is wrong with the answer I put for 1st response. Going to work on the Genetics homework before the 9 am lecture after making myself breakfast. Took measurements at 7:22 am, and same weight last 3 days but bloat in belly of waistline could be the cupcakes and gluten ate yesterday and/or the waist trimmer being set to 31" instead of 30" as it made me get indigestion a few days ago. Could also be my gut health from probiotics getting reduced from all the BMs the other day in indigestion. I don''t take probiotics, but the bacteria that we have in our guts is supposed to be there for gut health and probiotics is a way to keep it healthy is the current knowledge being passed around. It rained last night, wet outside and most likely not working out today. I have to get more stuff for the house. Did laundry and found $20 in the dryer, because I left a tip from one of my clients in my shirt or pants pocket a few days ago. I ordered some action figures for the roommate earlier and got the email saying they changed pick up time from today to the 24th so I then cancelled the order. His birthday is the 17th. And I ordered them because I could pick them up between now and on the 17th. Thats too bad. Will just get him a gift card, once they refund me my money. They charged my card first because I got the notificiation as soon as I selected the purchase button. Made 3 scrambled eggs with 1/3 green bell pepper and 2 tbs sourcream in a ninja blender then scrambled in 2 tbs olive oil with 4 mini corn tortilla quesadillas with mozzarella cheese. Need more cheese. Making my 3rd cup of coffee too, because starting to feel a headache coming on. I did all studying planned, and the lecture ran 1/2 hour past time with a beginning that kicked out the instructor and froze only her screen the first 10 minutes 2X, I ate 2 of the birthday cupcakes I got Shane and me, because its not his birthday yet not until next Wednesday and thats the frozen cake. He hasn''t even touched one. Taste just like the ones I ate yesterday for my coworker''s going away. I ate gluten and processed sweets today because of those cupcakes. I also made a smoothie with banana, serving of pea protein, honey 1 tbs, 1 serving walnuts about 13 pcs or 1/3 cup, and bottled water because I didn''t have almond milk. Growly didn''t eat any must have been the cocoa 2 tbsp and no almond milk. He must like the almond milk. He hasn''t been feeling well was under the bed the beginning of the day and when I woke up, he didn''t eat his breakfast, and I don''t think his lunch but did have some of the scrambled eggs and quesadilla I shared with him earlier. He probably isn''t feeling well but looks the same. Acts the same too. Went to get my car smogged for $40 locally and was in and out in 30 minutes tops and passed, uploaded it to Tred, and coincidentally today Tred said was the last day to upload it. It lasts 60 days for new owner if I sell it during that time period. It should pass, because I take great care of it, its new, and its still under warranty. But still great to know it does as expected and passes smog. Its a cute car and priced at 6% below dealers and has low miles that any first time or decent dude or girl with credit could get a car loan for through Tred. I keep dropping the price $400 every month after making the payment. Its selling right now for $18,300 with 74,000 miles. It actually has 100 miles less that 74k. Runs great. I have lecture in chemistry soon, and want to work on the lab and/or homework. No workout today. But the sun is out, and its cloudy but spotted cloudy. Supposed to continue intermittant showers throughout the day. Last time it rained was around 7 am and its almost 2 pm. I had 4 mini mozzarella corn tortilla quesadillas around 2 pm before lecture and about 2 servings of the Tostitos chips with about 10 tbs sourcream, then after lecture made a pot of spaghetti that makes about 5 bowls using 2 pkgs of Barilla gluten free spaghetti noodles, 1 pkg beyond meat, 1 jar prego 3 cheese sauce, 2 cups or about 15 white and sliced mushrooms, 2 bell peppers one yellow and one green, and 2 tbs olive oil. Had a bowl with 1/4 cup mozzarella and parmesan cheese blend from Target the Good & Gather brand shredded cheese. At break time cleaned Growly''s butt tail because he was dirty and did dishes before that, because I used the sink then wiped it down afterwards with disinfectant wipes and covered him in a towel after using a separate towel to dry him off. The moisture comes off him very well because he is a poodle. Finished class, got to midway or a third of the problems of lab part II on the unkown and confused by questions. Bed time by 10 pm.'

$wsMeasures.Range("AB57").Formula = "=210+60+240+13.3+160+267+270+150+60+80+130+105+20+60+200+7.5+266.67+80+614+100+300+30"
$wsMeasures.Range("AC57").Formula = "=15+5+28+0+10+2.67+12+7+5+5+2+0+1+0+20+0.625+2.67+5+24.8+6+14+25"
$wsMeasures.Range("AD57").Formula = "=4.5+3.5+4+0+7+0+2.6+1+3.5+3.5+0+0+0+0+2+0+0+3.5+6.1+4+2+17.5"
$wsMeasures.Range("AE57").Formula = "=18+1+0+0.3+12+5.33+3+2+1+6+18+1+2+0+5+0.25+5.33+6+26.1+8+4+5"
$wsMeasures.Range("AF57").Formula = "=0+2+0+3.3+2+53.33+36+18+2+1+9+27+6+17+4+0.25+53.33+1+72.9+2+36+10"
$wsMeasures.Range("AG57").Formula = "=0+0+0+1+0+5.33+0+1+0+0+2+3+2+0+2+0+5.33+0+5.5+0+2+0"
$wsMeasures.Range("AH57").Formula = "=210+15+0+0+380+53.33+310+130+15+190+320+1+0+0+0+28.75+53.33+190+798+280+260+75"

# --- sheet view / selection bookkeeping (cosmetic, matches authored commit) --
$wsMeasures.Application.ActiveWindow.ScrollColumn = 25  # topLeftCell Y1 (col 25)
$wsMeasures.Range("Z58").Select()
$wsNutrition.Range("B146:H146").Select()

$wb.Save()
